# The deck currently carries two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colour scheme (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral" colour scheme      (used by the slide master /
#                                                           the presentation's main design)
#
# The authored edit swaps the two themes' contents, so the design that drives the
# slides/slide master switches from "Integral" colours to the stock "Office Theme"
# colours (and vice versa for the notes master's theme).
#
# The presentation's single exposed Design/Theme (reached from any slide's
# ThemeColorScheme) is the one backed by ppt/theme/theme2.xml, so we recolour it,
# slot by slot, to the twelve "Office Theme" colours that used to live in
# ppt/theme/theme1.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477   # folHlink -> 954F72
